$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-02-23 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-02-24 Monday", 2)

# Update table cells by position (row, column) since some values repeat
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "83÷3=27, 2"
$t.Cell(1,2).Range.Text = "88÷7=12, 4"
$t.Cell(1,3).Range.Text = "46÷9=5, 1"
$t.Cell(1,4).Range.Text = "17÷7=2, 3"
$t.Cell(1,5).Range.Text = "67÷7=9, 4"

$t.Cell(5,1).Range.Text = "57÷9=6, 3"
$t.Cell(5,2).Range.Text = "75÷4=18, 3"
$t.Cell(5,3).Range.Text = "88÷8=11, 0"
$t.Cell(5,4).Range.Text = "20÷2=10, 0"
$t.Cell(5,5).Range.Text = "70÷9=7, 7"

$t.Cell(9,1).Range.Text = "69÷3=23, 0"
$t.Cell(9,2).Range.Text = "18÷3=6, 0"
$t.Cell(9,3).Range.Text = "94÷4=23, 2"
$t.Cell(9,4).Range.Text = "81÷9=9, 0"
$t.Cell(9,5).Range.Text = "67÷7=9, 4"

$t.Cell(13,1).Range.Text = "41÷4=10, 1"
$t.Cell(13,2).Range.Text = "95÷7=13, 4"
$t.Cell(13,3).Range.Text = "19÷4=4, 3"
$t.Cell(13,4).Range.Text = "39÷6=6, 3"
$t.Cell(13,5).Range.Text = "98÷8=12, 2"

$t.Cell(17,1).Range.Text = "83÷7=11, 6"
$t.Cell(17,2).Range.Text = "15÷9=1, 6"
$t.Cell(17,3).Range.Text = "67÷8=8, 3"
$t.Cell(17,4).Range.Text = "57÷3=19, 0"
$t.Cell(17,5).Range.Text = "91÷5=18, 1"
